$d = $word.ActiveDocument

# --- Location 1: Professional summary paragraph ---
$d.Content.Find.Execute(
    "Data engineering professional with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed geospatial ML algorithms improving classification accuracy from 23% to 64%. Built Civic Graph data warehouse processing billions of records and platforms serving thousands of analysts nationwide.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data engineering professional with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting 50M voters, developed geospatial ML algorithms improving classification accuracy from 23% to 64%. Built Civic Graph data warehouse processing billions of records and platforms serving thousands of analysts nationwide.",
    2)

# --- Location 2: Experience bullet - replace text, then re-format the "50M" run as bold + colored ---
$d.Content.Find.Execute(
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Discovered systematic race coding errors affecting 50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from ",
    2)

# Narrow in on the newly-inserted "50M" token within that bullet and give it the same
# bold + dark slate color formatting used for the other highlighted stats in this bullet.
$bulletRange = $d.Content
$bulletRange.Find.Execute("affecting 50M voters, developed geospatial machine learning")
$hit = $bulletRange.Duplicate
$hit.Find.Execute("50M")
$hit.Font.Bold = 1
$hit.Font.Color = 5258796

# --- Location 3: Key Projects "Impact:" line ---
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2)
